$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing data ends at row 1079 (date 2023-11-29) with fixed values that
# repeat for the new rows being appended (2023-11-30 .. 2023-12-20).
$lastRow = 1079
$startDate = Get-Date -Year 2023 -Month 11 -Day 30 -Hour 0 -Minute 0 -Second 0

$valB = 353806241
$valC = 60400000
$valD = 340000000
$valE = 263200000

for ($i = 0; $i -lt 21; $i++) {
    $row = $lastRow + 1 + $i
    $d = $startDate.AddDays($i)
    $dateStr = $d.ToString("yyyy-MM-ddTHH:mm:ssZ")

    $ws.Cells.Item($row, 1).Value = $dateStr
    $ws.Cells.Item($row, 2).Value = $valB
    $ws.Cells.Item($row, 3).Value = $valC
    $ws.Cells.Item($row, 4).Value = $valD
    $ws.Cells.Item($row, 5).Value = $valE
}
